# Update results for Steel
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B3: Hydrogen / Iron & steel demand value
$ws.Range("B3").Value = 389.7088409463599

# C4: Methanol / Chemicals demand value
$ws.Range("C4").Value = 30.49981016068242
